# Update Name of Algo
# Applies corrected imputed values to the KNN result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.786
$ws.Range("B4").Value = 5.340999999999999
$ws.Range("A6").Value = -22.209
$ws.Range("A7").Value = -20.296
$ws.Range("C7").Value = -11.879
$ws.Range("C8").Value = -12.283
$ws.Range("B9").Value = 5.808
$ws.Range("C10").Value = -13.173
$ws.Range("B12").Value = 4.939
$ws.Range("C13").Value = -13.097
$ws.Range("A16").Value = -21.2
$ws.Range("C16").Value = -12.883
$ws.Range("B17").Value = 5.644
$ws.Range("B18").Value = 5.988999999999999
$ws.Range("B19").Value = 7.309
$ws.Range("A20").Value = -21.731
$ws.Range("B20").Value = 5.773000000000001
$ws.Range("B26").Value = 5.681999999999999
$ws.Range("A28").Value = -21.619
$ws.Range("A29").Value = -21.5
$ws.Range("C30").Value = -12.664
$ws.Range("B31").Value = 5.965999999999999
$ws.Range("A32").Value = -21.528
$ws.Range("B39").Value = 6.639999999999999
$ws.Range("A40").Value = -21.432
$ws.Range("B40").Value = 6.167
$ws.Range("C40").Value = -11.547
$ws.Range("B41").Value = 6.281000000000001
$ws.Range("B42").Value = 6.256
$ws.Range("B43").Value = 6.071000000000001
$ws.Range("C44").Value = -12.168
$ws.Range("A46").Value = -21.59
$ws.Range("B47").Value = 5.455
$ws.Range("B48").Value = 5.443
$ws.Range("A51").Value = -21.557
$ws.Range("A52").Value = -21.651
$ws.Range("A57").Value = -21.889
$ws.Range("A59").Value = -22.119
$ws.Range("A62").Value = -21.779
$ws.Range("B63").Value = 5.252000000000001
$ws.Range("B64").Value = 5.527
$ws.Range("A66").Value = -21.563
$ws.Range("A73").Value = -21.021
$ws.Range("A74").Value = -20.504
$ws.Range("B76").Value = 6.399000000000001
$ws.Range("B81").Value = 5.859
$ws.Range("B89").Value = 5.315
$ws.Range("C89").Value = -14.055
$ws.Range("C91").Value = -12.783
$ws.Range("A92").Value = -21.566
$ws.Range("B94").Value = 5.795000000000001
$ws.Range("A100").Value = -21.885
